$wb = $excel.ActiveWorkbook

# Update the status string: "Ready for handoff" -> "Handed back: in sync with en-us"
# This shared string is used by B2 on both the zh-cn and de-de sheets.

$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$zhcn.Range("B2").Value = "Handed back: in sync with en-us"

$zhcn.Range("E2").Value = "6ad85588-6d3d-461e-8878-6a304486d3a7.md"
$zhcn.Range("F2").Value = "6ad85588-6d3d-461e-8878-6a304486d3a7.57ca24cbeb9804a84720f41594fce91a7dcad80c.zh-cn.xlf"
$zhcn.Range("G2").Value = "2016-01-09 03:39:25"
$zhcn.Range("H2").Value = "Include"

$zhcn.Hyperlinks.Add($zhcn.Range("E2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/x/6ad85588-6d3d-461e-8878-6a304486d3a7.md", "", "", "6ad85588-6d3d-461e-8878-6a304486d3a7.md")
$zhcn.Hyperlinks.Add($zhcn.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/x/6ad85588-6d3d-461e-8878-6a304486d3a7.57ca24cbeb9804a84720f41594fce91a7dcad80c.zh-cn.xlf", "", "", "6ad85588-6d3d-461e-8878-6a304486d3a7.57ca24cbeb9804a84720f41594fce91a7dcad80c.zh-cn.xlf")

$dede.Range("B2").Value = "Handed back: in sync with en-us"

$dede.Range("E2").Value = "6ad85588-6d3d-461e-8878-6a304486d3a7.md"
$dede.Range("F2").Value = "6ad85588-6d3d-461e-8878-6a304486d3a7.57ca24cbeb9804a84720f41594fce91a7dcad80c.de-de.xlf"
$dede.Range("G2").Value = "2016-01-09 03:39:42"
$dede.Range("H2").Value = "Include"

$dede.Hyperlinks.Add($dede.Range("E2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/x/6ad85588-6d3d-461e-8878-6a304486d3a7.md", "", "", "6ad85588-6d3d-461e-8878-6a304486d3a7.md")
$dede.Hyperlinks.Add($dede.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/x/6ad85588-6d3d-461e-8878-6a304486d3a7.57ca24cbeb9804a84720f41594fce91a7dcad80c.de-de.xlf", "", "", "6ad85588-6d3d-461e-8878-6a304486d3a7.57ca24cbeb9804a84720f41594fce91a7dcad80c.de-de.xlf")
